# Atualização do servidor e relatórios
# Adds a "Totais" (Totals) summary row above the existing header/data
# block on the "Registros" sheet, pushing the previous content down by
# two rows (one new label/value row + one blank spacer row), and pads
# the sheet with two more blank rows at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: insert two new blank rows above the current row 1.
#    Everything that used to live in rows 1-998 now lives in rows 3-1000,
#    carrying its original formatting down with it.
# ---------------------------------------------------------------------
$ws.Rows("1:2").Insert()

# Restore the custom row height on the two freshly-inserted rows (Excel
# does not automatically carry this over to brand new rows).
$ws.Rows("1:2").RowHeight = 15.75

# ---------------------------------------------------------------------
# 2) Populate the new "Totais" row (row 1) with labels/placeholders.
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Totais"
$ws.Range("C1").Value = "{d.totais.publicacoes}"
$ws.Range("D1").Value = "{d.totais.videos}"
$ws.Range("E1").Value = "{d.totais.horas}"
$ws.Range("F1").Value = "{d.totais.revisitas}"
$ws.Range("G1").Value = "{d.totais.estudos}"

# ---------------------------------------------------------------------
# 3) Formatting.
# ---------------------------------------------------------------------

# Whole band A1:H2 -> bold 14pt Calibri (theme text colour), solid white
# fill, no border - this is the base look for the new "Totais" band.
$band = $ws.Range("A1:H2")
$band.Font.Bold = $true
$band.Font.Size = 14
$band.Font.Name = "Calibri"
$band.Font.ThemeColor = 1
$band.Interior.Color = 16777215

# B1 ("Totais" label) -> same bold font but explicit black colour, light
# blue fill and a thin black box border, right aligned.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Font.Size = 14
$b1.Font.Name = "Calibri"
$b1.Font.Color = 0
$b1.Interior.Color = 16308937
$b1.HorizontalAlignment = -4152
$b1.BorderAround(1, 2, 1, 0)

# C1:G1 (totals placeholders) -> same bold black font, stays on the white
# fill, right aligned, no border.
$totalsRange = $ws.Range("C1:G1")
$totalsRange.Font.Bold = $true
$totalsRange.Font.Size = 14
$totalsRange.Font.Name = "Calibri"
$totalsRange.Font.Color = 0
$totalsRange.Interior.Color = 16777215
$totalsRange.HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 4) Two extra blank rows at the bottom of the sheet.
# ---------------------------------------------------------------------
$lastRow = $ws.Rows.Count
$ws.Range("A999").RowHeight = 15.75
$ws.Range("A1000").RowHeight = 15.75
